$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Update the RGB source values (columns A-C, rows 1-4).
# The dependent formulas in D:H recalculate automatically.
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 7
$ws.Range("C1").Value = 117

$ws.Range("A2").Value = 251
$ws.Range("B2").Value = 21
$ws.Range("C2").Value = 124

$ws.Range("A3").Value = 254
$ws.Range("B3").Value = 253
$ws.Range("C3").Value = 127

$ws.Range("A4").Value = 42
$ws.Range("B4").Value = 253
$ws.Range("C4").Value = 122

$excel.CalculateFullRebuild()
